$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 2.300632911392405
$ws.Range("P2").Value = 3.145387085179154
$ws.Range("R2").Value = 15
$ws.Range("S2").Value = 23
$ws.Range("T2").Value = 34
$ws.Range("U2").Value = 135
$ws.Range("V2").Value = 26.66455696202532
$ws.Range("W2").Value = 18.32999476423991
$ws.Range("Y2").Value = 27
$ws.Range("Z2").Value = 39.5
$ws.Range("AA2").Value = 53.25
$ws.Range("AB2").Value = 223
$ws.Range("AC2").Value = 49.1740506329114
$ws.Range("AD2").Value = 39.89666920159571
$ws.Range("AI2").Value = 16
$ws.Range("AJ2").Value = 3.164556962025316
$ws.Range("AK2").Value = 2.82249638349463
$ws.Range("AM2").Value = 0.3333333333333333
$ws.Range("AQ2").Value = 0.5307215824371054
$ws.Range("AR2").Value = 0.3479943612339447
$ws.Range("AS2").Value = 24.25
$ws.Range("AU2").Value = 47.5
$ws.Range("AV2").Value = 48.5
$ws.Range("AX2").Value = 47.07499209651108
$ws.Range("AY2").Value = 18.10781587709074
$ws.Range("V3").Value = 43.45569620253165
$ws.Range("W3").Value = 38.9706596897797
$ws.Range("Y3").Value = 47.75
$ws.Range("AC3").Value = 92.11075949367088
$ws.Range("AD3").Value = 84.25051690139057
$ws.Range("AJ3").Value = 4.193037974683544
$ws.Range("AK3").Value = 4.040754559247215
$ws.Range("K5").Value = 1
$ws.Range("N5").Value = 17
$ws.Range("O5").Value = 3.522151898734177
$ws.Range("P5").Value = 4.576282866961004
$ws.Range("T5").Value = 29.25
$ws.Range("U5").Value = 98
$ws.Range("V5").Value = 21.92721518987342
$ws.Range("W5").Value = 15.29262982038693
$ws.Range("Y5").Value = 27
$ws.Range("Z5").Value = 38
$ws.Range("AA5").Value = 54
$ws.Range("AB5").Value = 190
$ws.Range("AC5").Value = 46.2626582278481
$ws.Range("AD5").Value = 32.5537016336737
$ws.Range("AF5").Value = 1
$ws.Range("AH5").Value = 3.25
$ws.Range("AJ5").Value = 3.091772151898734
$ws.Range("AK5").Value = 2.610320546876234
$ws.Range("AM5").Value = 0.6666666666666666
$ws.Range("AQ5").Value = 0.8365335807874116
$ws.Range("AR5").Value = 0.251086959137955
$ws.Range("AS5").Value = 23.77777777777778
$ws.Range("AX5").Value = 33.03482243432914
$ws.Range("AY5").Value = 6.477471883527956
$ws.Range("O6").Value = 7.186708860759493
$ws.Range("P6").Value = 8.515901871767431
$ws.Range("T6").Value = 13
$ws.Range("V6").Value = 9.882911392405063
$ws.Range("W6").Value = 9.952054541255478
$ws.Range("AA6").Value = 17
$ws.Range("AC6").Value = 9.113924050632912
$ws.Range("AD6").Value = 9.847713252272086
$ws.Range("AH6").Value = 2
$ws.Range("AJ6").Value = 0.9240506329113924
$ws.Range("AK6").Value = 0.9729307960672774
$ws.Range("AM6").Value = 1.473684210526316
$ws.Range("AP6").Value = 4
$ws.Range("AQ6").Value = 1.881864471304844
$ws.Range("AR6").Value = 0.6722739449024928
$ws.Range("AU6").Value = 18.60869565217391
$ws.Range("AV6").Value = 20.08333333333333
$ws.Range("AX6").Value = 19.04728065130852
$ws.Range("AY6").Value = 5.882771635131574
$ws.Range("O7").Value = 2.738589211618257
$ws.Range("P7").Value = 3.207108718880054
$ws.Range("R7").Value = 16
$ws.Range("V7").Value = 25.85477178423237
$ws.Range("W7").Value = 13.93645056025866
$ws.Range("Z7").Value = 37
$ws.Range("AC7").Value = 39.8298755186722
$ws.Range("AD7").Value = 21.95108434074362
$ws.Range("AJ7").Value = 2.842323651452282
$ws.Range("AK7").Value = 1.655808335772374
$ws.Range("AM7").Value = 0.4285714285714285
$ws.Range("AO7").Value = 0.6666666666666666
$ws.Range("AQ7").Value = 0.5773406545409167
$ws.Range("AR7").Value = 0.2964972793141106
$ws.Range("AU7").Value = 47.5
$ws.Range("AX7").Value = 44.44608321600023
$ws.Range("AY7").Value = 15.39743138309131
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 38.5
$ws.Range("V8").Value = 49.15644171779141
$ws.Range("W8").Value = 63.67318047830457
$ws.Range("Y8").Value = 0
$ws.Range("AC8").Value = 83.8159509202454
$ws.Range("AD8").Value = 101.0177890968519
$ws.Range("AF8").Value = 0
$ws.Range("AJ8").Value = 3.187116564417178
$ws.Range("AK8").Value = 2.988998231757231
$ws.Range("K10").Value = 1
$ws.Range("O10").Value = 3.402489626556016
$ws.Range("P10").Value = 3.883704474158696
$ws.Range("V10").Value = 23.44398340248962
$ws.Range("W10").Value = 10.74311674514074
$ws.Range("AA10").Value = 48
$ws.Range("AC10").Value = 37.43568464730291
$ws.Range("AD10").Value = 20.12867493486523
$ws.Range("AJ10").Value = 2.692946058091287
$ws.Range("AK10").Value = 1.601142831827155
$ws.Range("AN10").Value = 0.7368421052631579
$ws.Range("AQ10").Value = 0.7439813226011086
$ws.Range("AR10").Value = 0.2806086267855832
$ws.Range("AV10").Value = 42.8
$ws.Range("AX10").Value = 37.10115498771832
$ws.Range("AY10").Value = 11.126796191337
$ws.Range("N11").Value = 32
$ws.Range("O11").Value = 6.950207468879668
$ws.Range("P11").Value = 5.845725820926259
$ws.Range("V11").Value = 15.09128630705394
$ws.Range("W11").Value = 14.51952589062467
$ws.Range("AB11").Value = 44
$ws.Range("AC11").Value = 10.46058091286307
$ws.Range("AD11").Value = 9.324670574760313
$ws.Range("AI11").Value = 5
$ws.Range("AJ11").Value = 1.095435684647303
$ws.Range("AK11").Value = 0.9721217757503259
$ws.Range("AQ11").Value = 1.917173639492923
$ws.Range("AR11").Value = 0.8280222081088773
$ws.Range("AT11").Value = 14.66666666666667
$ws.Range("AU11").Value = 18.6
$ws.Range("AX11").Value = 19.95054735876594
$ws.Range("AY11").Value = 9.828944512850548
